$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 397
$ws.Range("I2").Value = 397
$ws.Range("K2").Value = 397
$ws.Range("M2").Value = -284

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H38").Value = 1630.1428
$ws.Range("I38").Value = 424.66666
$ws.Range("J38").Value = 3800
$ws.Range("K38").Value = 1273.99998
$ws.Range("L38").Value = 11400
$ws.Range("M38").Value = -901.9999800000001
$ws.Range("N38").Value = -12144

$ws.Range("H51").Value = 2900.3333
$ws.Range("J51").Value = 2900.3333
$ws.Range("L51").Value = 2900.3333
$ws.Range("N51").Value = -3868.3333

$ws.Range("H58").Value = 2536.7273
$ws.Range("I58").Value = 307.5
$ws.Range("J58").Value = 2759.65
$ws.Range("K58").Value = 922.5
$ws.Range("L58").Value = 8278.950000000001
$ws.Range("M58").Value = -772.5
$ws.Range("N58").Value = -8578.950000000001

$ws.Range("H87").Value = 30900
$ws.Range("J87").Value = 30900
$ws.Range("L87").Value = 30900
$ws.Range("N87").Value = -33396

$ws.Range("H90").Value = 30900
$ws.Range("J90").Value = 30900
$ws.Range("L90").Value = 92700
$ws.Range("N90").Value = -105180

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H111").Value = 667.2143
$ws.Range("I111").Value = 587.7
$ws.Range("K111").Value = 1763.1
$ws.Range("M111").Value = 1303.9

$ws.Range("H116").Value = 10281.883
$ws.Range("I116").Value = 15334.444
$ws.Range("J116").Value = 4597.75
$ws.Range("K116").Value = 15334.444
$ws.Range("L116").Value = 4597.75
$ws.Range("M116").Value = -11892.444
$ws.Range("N116").Value = -11481.75

$ws.Range("H138").Value = 1954.2084
$ws.Range("I138").Value = 1163.75
$ws.Range("J138").Value = 2518.8215
$ws.Range("K138").Value = 3491.25
$ws.Range("L138").Value = 7556.4645
$ws.Range("M138").Value = 1648.75
$ws.Range("N138").Value = -17836.4645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4316.6665
$ws.Range("I3").Value = 2225
$ws.Range("K3").Value = 2225
$ws.Range("M3").Value = -2110

$ws.Range("H61").Value = 186010.45
$ws.Range("I61").Value = 4802.0884
$ws.Range("J61").Value = 479395.44
$ws.Range("K61").Value = 4802.0884
$ws.Range("L61").Value = 479395.44
$ws.Range("M61").Value = -4590.0884
$ws.Range("N61").Value = -479819.44

$ws.Range("H63").Value = 200002000
$ws.Range("I63").Value = 250002000
$ws.Range("J63").Value = 2006
$ws.Range("K63").Value = 250002000
$ws.Range("L63").Value = 2006
$ws.Range("M63").Value = -250001314
$ws.Range("N63").Value = -3378

$ws.Range("H66").Value = 200002000
$ws.Range("I66").Value = 250002000
$ws.Range("J66").Value = 2006
$ws.Range("K66").Value = 1250010000
$ws.Range("L66").Value = 10030
$ws.Range("M66").Value = -1250006568
$ws.Range("N66").Value = -16894

$ws.Range("H136").Value = 186010.45
$ws.Range("I136").Value = 4802.0884
$ws.Range("J136").Value = 479395.44
$ws.Range("K136").Value = 14406.2652
$ws.Range("L136").Value = 1438186.32
$ws.Range("M136").Value = -11856.2652
$ws.Range("N136").Value = -1443286.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1401.8572
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H18").Value = 7210
$ws.Range("J18").Value = 7210
$ws.Range("L18").Value = 7210
$ws.Range("N18").Value = -8268

$ws.Range("H24").Value = 2744
$ws.Range("I24").Value = 2744
$ws.Range("K24").Value = 2744
$ws.Range("M24").Value = -2509

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H64").Value = 461
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 461
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 461
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -911

$ws.Range("H67").Value = 461
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 461
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 461
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2021

$ws.Range("H75").Value = 20771.285
$ws.Range("I75").Value = 10475.667
$ws.Range("J75").Value = 28493
$ws.Range("K75").Value = 10475.667
$ws.Range("L75").Value = 28493
$ws.Range("M75").Value = -9539.666999999999
$ws.Range("N75").Value = -30365

$ws.Range("H78").Value = 20771.285
$ws.Range("I78").Value = 10475.667
$ws.Range("J78").Value = 28493
$ws.Range("K78").Value = 31427.001
$ws.Range("L78").Value = 85479
$ws.Range("M78").Value = -26747.001
$ws.Range("N78").Value = -94839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 300
$ws.Range("K15").Value = 300
$ws.Range("M15").Value = -130

$ws.Range("H16").Value = 2580.875
$ws.Range("I16").Value = 2340.3333
$ws.Range("J16").Value = 2725.2
$ws.Range("K16").Value = 2340.3333
$ws.Range("L16").Value = 2725.2
$ws.Range("M16").Value = -2053.3333
$ws.Range("N16").Value = -3299.2

$ws.Range("H31").Value = 4558.0645
$ws.Range("I31").Value = 1703.6923
$ws.Range("J31").Value = 6619.5557
$ws.Range("K31").Value = 1703.6923
$ws.Range("L31").Value = 6619.5557
$ws.Range("M31").Value = -1408.6923
$ws.Range("N31").Value = -7209.5557

$ws.Range("H34").Value = 4558.0645
$ws.Range("I34").Value = 1703.6923
$ws.Range("J34").Value = 6619.5557
$ws.Range("K34").Value = 1703.6923
$ws.Range("L34").Value = 6619.5557
$ws.Range("M34").Value = -1501.6923
$ws.Range("N34").Value = -7023.5557

$ws.Range("H62").Value = 6421.5264
$ws.Range("I62").Value = 6833.1665
$ws.Range("J62").Value = 5715.857
$ws.Range("K62").Value = 6833.1665
$ws.Range("L62").Value = 5715.857
$ws.Range("M62").Value = -6209.1665
$ws.Range("N62").Value = -6963.857

$ws.Range("H65").Value = 6421.5264
$ws.Range("I65").Value = 6833.1665
$ws.Range("J65").Value = 5715.857
$ws.Range("K65").Value = 34165.8325
$ws.Range("L65").Value = 28579.285
$ws.Range("M65").Value = -31045.8325
$ws.Range("N65").Value = -34819.285

$ws.Range("H113").Value = 2580.875
$ws.Range("I113").Value = 2340.3333
$ws.Range("J113").Value = 2725.2
$ws.Range("K113").Value = 2340.3333
$ws.Range("L113").Value = 2725.2
$ws.Range("M113").Value = -170.3332999999998
$ws.Range("N113").Value = -7065.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4519.7744
$ws.Range("I5").Value = 6987.7334
$ws.Range("J5").Value = 2206.0625
$ws.Range("K5").Value = 20963.2002
$ws.Range("L5").Value = 6618.1875
$ws.Range("M5").Value = -20851.2002
$ws.Range("N5").Value = -6842.1875

$ws.Range("H121").Value = 20800
$ws.Range("J121").Value = 82306.5
$ws.Range("L121").Value = 246919.5
$ws.Range("N121").Value = -249539.5

$ws.Range("H122").Value = 4820.08
$ws.Range("I122").Value = 750.8889
$ws.Range("J122").Value = 7109
$ws.Range("K122").Value = 6758.0001
$ws.Range("L122").Value = 63981
$ws.Range("M122").Value = -4308.0001
$ws.Range("N122").Value = -68881

$ws.Range("H130").Value = 6971.838
$ws.Range("I130").Value = 2850
$ws.Range("K130").Value = 8550
$ws.Range("M130").Value = -3530

$ws.Range("H135").Value = 4519.7744
$ws.Range("I135").Value = 6987.7334
$ws.Range("J135").Value = 2206.0625
$ws.Range("K135").Value = 62889.6006
$ws.Range("L135").Value = 19854.5625
$ws.Range("M135").Value = -60354.6006
$ws.Range("N135").Value = -24924.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10002807
$ws.Range("I20").Value = 15000000
$ws.Range("K20").Value = 15000000
$ws.Range("M20").Value = -14999755

$ws.Range("H24").Value = 3338666.8

$ws.Range("H113").Value = 34484030
$ws.Range("I113").Value = 41667640
$ws.Range("J113").Value = 2718
$ws.Range("K113").Value = 41667640
$ws.Range("L113").Value = 2718
$ws.Range("M113").Value = -41665470
$ws.Range("N113").Value = -7058

$ws.Range("H122").Value = 40957556
$ws.Range("I122").Value = 62637830
$ws.Range("J122").Value = 5921.6665
$ws.Range("K122").Value = 187913490
$ws.Range("L122").Value = 17764.9995
$ws.Range("M122").Value = -187911040
$ws.Range("N122").Value = -22664.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2600.818
$ws.Range("I7").Value = 2143.4285
$ws.Range("J7").Value = 3401.25
$ws.Range("K7").Value = 2143.4285
$ws.Range("L7").Value = 3401.25
$ws.Range("M7").Value = -2031.4285
$ws.Range("N7").Value = -3625.25

$ws.Range("H122").Value = 5435566
$ws.Range("I122").Value = 5960457.5
$ws.Range("K122").Value = 17881372.5
$ws.Range("M122").Value = -17878922.5

$ws.Range("H126").Value = 2600.818
$ws.Range("I126").Value = 2143.4285
$ws.Range("J126").Value = 3401.25
$ws.Range("K126").Value = 6430.2855
$ws.Range("L126").Value = 10203.75
$ws.Range("M126").Value = -3960.2855
$ws.Range("N126").Value = -15143.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 10000
$ws.Range("J31").Value = 10000
$ws.Range("L31").Value = 10000
$ws.Range("N31").Value = -10696

$ws.Range("H100").Value = 338
$ws.Range("I100").Value = 338
$ws.Range("K100").Value = 676
$ws.Range("M100").Value = -135

$ws.Range("H107").Value = 47619796
$ws.Range("I107").Value = 90909680
$ws.Range("J107").Value = 919.3
$ws.Range("K107").Value = 272729040
$ws.Range("L107").Value = 2757.9
$ws.Range("M107").Value = -272727120
$ws.Range("N107").Value = -6597.9

$ws.Range("H136").Value = 2349.644
$ws.Range("I136").Value = 2563.3547
$ws.Range("J136").Value = 2113.0356
$ws.Range("K136").Value = 7690.0641
$ws.Range("L136").Value = 6339.1068
$ws.Range("M136").Value = -5140.0641
$ws.Range("N136").Value = -11439.1068
